$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column G ("Recorded By") holds comma-separated lists of recorders.
# For every row whose value contains more than one entry, rotate the
# list left by one (move the first entry to the end).
$rows = 2..153

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Length -gt 1) {
            $newParts = $parts[1..($parts.Length - 1)] + $parts[0]
            $cell.Value = [string]::Join(", ", $newParts)
        }
    }
}
